$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 11 and 12, shifting existing data (rows 11+) down by two.
$ws.Rows("11:12").Insert()

# New row 11 data
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value = 45163
$ws.Cells.Item(11, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11, 5).Value = 15
$ws.Cells.Item(11, 6).Value = 100112031
$ws.Cells.Item(11, 7).Value = "Poroto verde"
$ws.Cells.Item(11, 8).Value = "Magnum"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 2200
$ws.Cells.Item(11, 11).Value = 1000
$ws.Cells.Item(11, 12).Value = 1200
$ws.Cells.Item(11, 13).Value = 1073
$ws.Cells.Item(11, 14).Value = "`$/kilo"
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 16).Value = 1073
$ws.Cells.Item(11, 17).Value = 1
$ws.Cells.Item(11, 18).Value = "Hortaliza"

# New row 12 data
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(12, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(12, 4).Value = 45163
$ws.Cells.Item(12, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = 100112031
$ws.Cells.Item(12, 7).Value = "Poroto verde"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 2300
$ws.Cells.Item(12, 11).Value = 1200
$ws.Cells.Item(12, 12).Value = 1300
$ws.Cells.Item(12, 13).Value = 1265
$ws.Cells.Item(12, 14).Value = "`$/kilo"
$ws.Cells.Item(12, 15).Value = "Perú"
$ws.Cells.Item(12, 16).Value = 1265
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = "Hortaliza"
